# Auto-generated edit script applying the cryptos.xlsx price/volume refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "92.197.62"; ForceText = $false },
    @{ Cell = "E2"; Value = "  -3.32%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "3.266.67"; ForceText = $false },
    @{ Cell = "E3"; Value = "  -5.56%  "; ForceText = $false },
    @{ Cell = "E4"; Value = "  +0.14%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "224.44"; ForceText = $true },
    @{ Cell = "E5"; Value = "  -6.89%  "; ForceText = $false },
    @{ Cell = "D6"; Value = "601.72"; ForceText = $true },
    @{ Cell = "E6"; Value = "  -6.68%  "; ForceText = $false },
    @{ Cell = "D7"; Value = "1.34"; ForceText = $true },
    @{ Cell = "E7"; Value = "  -8.98%  "; ForceText = $false },
    @{ Cell = "D8"; Value = "0.372"; ForceText = $true },
    @{ Cell = "E8"; Value = "  -7.75%  "; ForceText = $false },
    @{ Cell = "E9"; Value = "  +0.07%  "; ForceText = $false },
    @{ Cell = "D10"; Value = "0.898"; ForceText = $true },
    @{ Cell = "E10"; Value = "  -9.92%  "; ForceText = $false },
    @{ Cell = "D11"; Value = "3.262.86"; ForceText = $false },
    @{ Cell = "E11"; Value = "  -5.57%  "; ForceText = $false },
    @{ Cell = "D12"; Value = "40.91"; ForceText = $true },
    @{ Cell = "E12"; Value = "  -1.80%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "0.189"; ForceText = $true },
    @{ Cell = "E13"; Value = "  -4.56%  "; ForceText = $false },
    @{ Cell = "B14"; Value = "WrappedBTC"; ForceText = $false },
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; ForceText = $false },
    @{ Cell = "D14"; Value = "92.105.38"; ForceText = $false },
    @{ Cell = "E14"; Value = "  -3.23%  "; ForceText = $false },
    @{ Cell = "B15"; Value = "Toncoin"; ForceText = $false },
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; ForceText = $false },
    @{ Cell = "D15"; Value = "5.83"; ForceText = $true },
    @{ Cell = "E15"; Value = "  -4.72%  "; ForceText = $false },
    @{ Cell = "D16"; Value = "3.873.17"; ForceText = $false },
    @{ Cell = "E16"; Value = "  -5.80%  "; ForceText = $false },
    @{ Cell = "D17"; Value = "0.0000236"; ForceText = $true },
    @{ Cell = "E17"; Value = "  -8.09%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "7.82"; ForceText = $true },
    @{ Cell = "E18"; Value = "  -7.49%  "; ForceText = $false },
    @{ Cell = "D19"; Value = "3.284.37"; ForceText = $false },
    @{ Cell = "E19"; Value = "  -5.11%  "; ForceText = $false },
    @{ Cell = "D20"; Value = "16.76"; ForceText = $true },
    @{ Cell = "E20"; Value = "  -6.39%  "; ForceText = $false },
    @{ Cell = "D21"; Value = "10.57"; ForceText = $true },
    @{ Cell = "E21"; Value = "  -8.07%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "3.36"; ForceText = $true },
    @{ Cell = "E22"; Value = "  +5.55%  "; ForceText = $false },
    @{ Cell = "D23"; Value = "480.55"; ForceText = $true },
    @{ Cell = "E23"; Value = "  -4.42%  "; ForceText = $false },
    @{ Cell = "D24"; Value = "0.431"; ForceText = $true },
    @{ Cell = "E24"; Value = "  -15.57%  "; ForceText = $false },
    @{ Cell = "D25"; Value = "0.0000174"; ForceText = $true },
    @{ Cell = "E25"; Value = "  -9.27%  "; ForceText = $false },
    @{ Cell = "D26"; Value = "5.92"; ForceText = $true },
    @{ Cell = "E26"; Value = "  -10.49%  "; ForceText = $false },
    @{ Cell = "D27"; Value = "88.30"; ForceText = $true },
    @{ Cell = "E27"; Value = "  -7.12%  "; ForceText = $false },
    @{ Cell = "D28"; Value = "3.450.47"; ForceText = $false },
    @{ Cell = "E28"; Value = "  -5.39%  "; ForceText = $false },
    @{ Cell = "D29"; Value = "11.40"; ForceText = $true },
    @{ Cell = "E29"; Value = "  -5.42%  "; ForceText = $false },
    @{ Cell = "E30"; Value = "  -0.06%  "; ForceText = $false },
    @{ Cell = "D31"; Value = "10.73"; ForceText = $true },
    @{ Cell = "E31"; Value = "  -8.37%  "; ForceText = $false },
    @{ Cell = "D32"; Value = "0.135"; ForceText = $true },
    @{ Cell = "E32"; Value = "  -1.43%  "; ForceText = $false },
    @{ Cell = "D33"; Value = "2.55"; ForceText = $true },
    @{ Cell = "E33"; Value = "  -7.61%  "; ForceText = $false },
    @{ Cell = "D34"; Value = "0.998"; ForceText = $true },
    @{ Cell = "E34"; Value = "  -0.23%  "; ForceText = $false },
    @{ Cell = "D35"; Value = "0.169"; ForceText = $true },
    @{ Cell = "E35"; Value = "  -8.28%  "; ForceText = $false },
    @{ Cell = "D36"; Value = "27.51"; ForceText = $true },
    @{ Cell = "E36"; Value = "  -11.64%  "; ForceText = $false },
    @{ Cell = "D37"; Value = "0.513"; ForceText = $true },
    @{ Cell = "E37"; Value = "  -9.90%  "; ForceText = $false },
    @{ Cell = "D38"; Value = "522.17"; ForceText = $true },
    @{ Cell = "E38"; Value = "  -1.01%  "; ForceText = $false },
    @{ Cell = "E39"; Value = "  -0.10%  "; ForceText = $false },
    @{ Cell = "D40"; Value = "7.13"; ForceText = $true },
    @{ Cell = "E40"; Value = "  -8.31%  "; ForceText = $false },
    @{ Cell = "D41"; Value = "0.144"; ForceText = $true },
    @{ Cell = "E41"; Value = "  -4.14%  "; ForceText = $false },
    @{ Cell = "D42"; Value = "1.32"; ForceText = $true },
    @{ Cell = "E42"; Value = "  -8.66%  "; ForceText = $false },
    @{ Cell = "B43"; Value = "WhiteBITCoin"; ForceText = $false },
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"; ForceText = $false },
    @{ Cell = "D43"; Value = "23.91"; ForceText = $true },
    @{ Cell = "E43"; Value = "  -0.85%  "; ForceText = $false },
    @{ Cell = "B44"; Value = "ARBITRUM"; ForceText = $false },
    @{ Cell = "C44"; Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; ForceText = $false },
    @{ Cell = "D44"; Value = "0.833"; ForceText = $true },
    @{ Cell = "E44"; Value = "  -8.68%  "; ForceText = $false },
    @{ Cell = "D45"; Value = "3.54"; ForceText = $true },
    @{ Cell = "E45"; Value = "  +1.44%  "; ForceText = $false },
    @{ Cell = "D46"; Value = "1.62"; ForceText = $true },
    @{ Cell = "E46"; Value = "  -4.93%  "; ForceText = $false },
    @{ Cell = "D47"; Value = "0.0400"; ForceText = $true },
    @{ Cell = "E47"; Value = "  -4.09%  "; ForceText = $false },
    @{ Cell = "D48"; Value = "5.20"; ForceText = $true },
    @{ Cell = "E48"; Value = "  -7.99%  "; ForceText = $false },
    @{ Cell = "D49"; Value = "51.54"; ForceText = $true },
    @{ Cell = "E49"; Value = "  -3.87%  "; ForceText = $false },
    @{ Cell = "D50"; Value = "2.04"; ForceText = $true },
    @{ Cell = "E50"; Value = "  -5.32%  "; ForceText = $false },
    @{ Cell = "D51"; Value = "7.68"; ForceText = $true },
    @{ Cell = "E51"; Value = "  -4.35%  "; ForceText = $false }
)

foreach ($item in $updates) {
    $c = $ws.Range($item.Cell)
    if ($item.ForceText) {
        $origStyle = $c.Style
        $c.NumberFormat = "@"
        $c.Value = $item.Value
        $c.Style = $origStyle
    } else {
        $c.Value = $item.Value
    }
}
